$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value2 = 3002.6924
$ws.Range("I39").Value2 = 950.7778
$ws.Range("K39").Value2 = 2852.3334
$ws.Range("M39").Value2 = -2556.3334

$ws.Range("H40").Value2 = 2389
$ws.Range("I40").Value2 = 2389
$ws.Range("J40").Value2 = 0
$ws.Range("K40").Value2 = 2389
$ws.Range("L40").Value2 = 0
$ws.Range("M40").Value2 = -2214
$ws.Range("N40").ClearContents()

$ws.Range("H51").Value2 = 12874.9
$ws.Range("I51").Value2 = 4899
$ws.Range("J51").Value2 = 13761.111
$ws.Range("K51").Value2 = 4899
$ws.Range("L51").Value2 = 13761.111
$ws.Range("M51").Value2 = -4415
$ws.Range("N51").Value2 = -14729.111

$ws.Range("H70").Value2 = 145840.14
$ws.Range("J70").Value2 = 334799.66
$ws.Range("L70").Value2 = 1004398.98
$ws.Range("N70").Value2 = -1004938.98

$ws.Range("H73").Value2 = 145840.14
$ws.Range("J73").Value2 = 334799.66
$ws.Range("L73").Value2 = 1004398.98
$ws.Range("N73").Value2 = -1006270.98

$ws.Range("H103").Value2 = 396.6
$ws.Range("I103").Value2 = 300
$ws.Range("J103").Value2 = 541.5
$ws.Range("K103").Value2 = 900
$ws.Range("L103").Value2 = 1624.5
$ws.Range("M103").Value2 = -314
$ws.Range("N103").Value2 = -2796.5

$ws.Range("H116").Value2 = 5908
$ws.Range("J116").Value2 = 7064.5713
$ws.Range("L116").Value2 = 7064.5713
$ws.Range("N116").Value2 = -13948.5713

$ws.Range("H137").Value2 = 0
$ws.Range("I137").Value2 = 0
$ws.Range("K137").Value2 = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 4683.972
$ws.Range("I61").Value2 = 3804.2334
$ws.Range("J61").Value2 = 9082.666999999999
$ws.Range("K61").Value2 = 3804.2334
$ws.Range("L61").Value2 = 9082.666999999999
$ws.Range("M61").Value2 = -3592.2334
$ws.Range("N61").Value2 = -9506.666999999999

$ws.Range("H63").Value2 = 2726.8
$ws.Range("I63").Value2 = 2819.5557
$ws.Range("J63").Value2 = 1892
$ws.Range("K63").Value2 = 2819.5557
$ws.Range("L63").Value2 = 1892
$ws.Range("M63").Value2 = -2133.5557
$ws.Range("N63").Value2 = -3264

$ws.Range("H66").Value2 = 2726.8
$ws.Range("I66").Value2 = 2819.5557
$ws.Range("J66").Value2 = 1892
$ws.Range("K66").Value2 = 14097.7785
$ws.Range("L66").Value2 = 9460
$ws.Range("M66").Value2 = -10665.7785
$ws.Range("N66").Value2 = -16324

$ws.Range("H74").Value2 = 4271.7144
$ws.Range("I74").Value2 = 3664.0908
$ws.Range("J74").Value2 = 6499.6665
$ws.Range("K74").Value2 = 3664.0908
$ws.Range("L74").Value2 = 6499.6665
$ws.Range("M74").Value2 = -2790.0908
$ws.Range("N74").Value2 = -8247.666499999999

$ws.Range("H77").Value2 = 4271.7144
$ws.Range("I77").Value2 = 3664.0908
$ws.Range("J77").Value2 = 6499.6665
$ws.Range("K77").Value2 = 18320.454
$ws.Range("L77").Value2 = 32498.3325
$ws.Range("M77").Value2 = -13952.454
$ws.Range("N77").Value2 = -41234.3325

$ws.Range("H122").Value2 = 4886.1387
$ws.Range("I122").Value2 = 4258.931
$ws.Range("K122").Value2 = 12776.793
$ws.Range("M122").Value2 = -10326.793

$ws.Range("H136").Value2 = 4683.972
$ws.Range("I136").Value2 = 3804.2334
$ws.Range("J136").Value2 = 9082.666999999999
$ws.Range("K136").Value2 = 11412.7002
$ws.Range("L136").Value2 = 27248.001
$ws.Range("M136").Value2 = -8862.700199999999
$ws.Range("N136").Value2 = -32348.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value2 = 458.8
$ws.Range("I22").Value2 = 458.8
$ws.Range("K22").Value2 = 458.8
$ws.Range("M22").Value2 = -285.8

$ws.Range("H58").Value2 = 27709
$ws.Range("J58").Value2 = 0
$ws.Range("L58").Value2 = 0
$ws.Range("N58").ClearContents()

$ws.Range("H105").Value2 = 3222.6667
$ws.Range("I105").Value2 = 2863.923
$ws.Range("K105").Value2 = 2863.923
$ws.Range("M105").Value2 = -1116.923

$ws.Range("H107").Value2 = 2368.5518
$ws.Range("I107").Value2 = 2058.72
$ws.Range("J107").Value2 = 4305
$ws.Range("K107").Value2 = 2058.72
$ws.Range("L107").Value2 = 4305
$ws.Range("M107").Value2 = -138.7199999999998
$ws.Range("N107").Value2 = -8145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 6816.05
$ws.Range("I31").Value2 = 7680.0713
$ws.Range("J31").Value2 = 6350.8076
$ws.Range("K31").Value2 = 7680.0713
$ws.Range("L31").Value2 = 6350.8076
$ws.Range("M31").Value2 = -7385.0713
$ws.Range("N31").Value2 = -6940.8076

$ws.Range("H34").Value2 = 6816.05
$ws.Range("I34").Value2 = 7680.0713
$ws.Range("J34").Value2 = 6350.8076
$ws.Range("K34").Value2 = 7680.0713
$ws.Range("L34").Value2 = 6350.8076
$ws.Range("M34").Value2 = -7478.0713
$ws.Range("N34").Value2 = -6754.8076

$ws.Range("H134").Value2 = 1116.1177
$ws.Range("I134").Value2 = 1031.6154
$ws.Range("K134").Value2 = 3094.8462
$ws.Range("M134").Value2 = -559.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 212338080
$ws.Range("I4").Value2 = 333563460
$ws.Range("K4").Value2 = 1000690380
$ws.Range("M4").Value2 = -1000690268

$ws.Range("H59").Value2 = 3232.8333
$ws.Range("I59").Value2 = 1199
$ws.Range("K59").Value2 = 3597
$ws.Range("M59").Value2 = -3057

$ws.Range("H80").Value2 = 4894
$ws.Range("I80").Value2 = 4900
$ws.Range("J80").Value2 = 4888
$ws.Range("K80").Value2 = 14700
$ws.Range("L80").Value2 = 14664
$ws.Range("M80").Value2 = -13764
$ws.Range("N80").Value2 = -16536

$ws.Range("H83").Value2 = 4894
$ws.Range("I83").Value2 = 4900
$ws.Range("J83").Value2 = 4888
$ws.Range("K83").Value2 = 44100
$ws.Range("L83").Value2 = 43992
$ws.Range("M83").Value2 = -39420
$ws.Range("N83").Value2 = -53352

$ws.Range("H113").Value2 = 2075.423
$ws.Range("J113").Value2 = 2216.4546
$ws.Range("L113").Value2 = 6649.3638
$ws.Range("N113").Value2 = -10989.3638

$ws.Range("H131").Value2 = 2176
$ws.Range("J131").Value2 = 2665.6667
$ws.Range("L131").Value2 = 7997.000100000001
$ws.Range("N131").Value2 = -18077.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value2 = 6464.968
$ws.Range("I126").Value2 = 5410.1816
$ws.Range("K126").Value2 = 16230.5448
$ws.Range("M126").Value2 = -13760.5448

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 905.1111
$ws.Range("I16").Value2 = 320.93332
$ws.Range("J16").Value2 = 3826
$ws.Range("K16").Value2 = 320.93332
$ws.Range("L16").Value2 = 3826
$ws.Range("M16").Value2 = -150.93332
$ws.Range("N16").Value2 = -4166

$ws.Range("H40").Value2 = 5810.1333
$ws.Range("I40").Value2 = 5835.885
$ws.Range("J40").Value2 = 5642.75
$ws.Range("K40").Value2 = 5835.885
$ws.Range("L40").Value2 = 5642.75
$ws.Range("M40").Value2 = -5699.885
$ws.Range("N40").Value2 = -5914.75

$ws.Range("H46").Value2 = 5100
$ws.Range("I46").Value2 = 0
$ws.Range("J46").Value2 = 5100
$ws.Range("K46").Value2 = 0
$ws.Range("L46").Value2 = 5100
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value2 = -5476

$ws.Range("H93").Value2 = 1388.08
$ws.Range("J93").Value2 = 1497.5
$ws.Range("L93").Value2 = 1497.5
$ws.Range("N93").Value2 = -3993.5

$ws.Range("H122").Value2 = 3373
$ws.Range("I122").Value2 = 1998
$ws.Range("J122").Value2 = 3569.4285
$ws.Range("K122").Value2 = 5994
$ws.Range("L122").Value2 = 10708.2855
$ws.Range("M122").Value2 = -3544
$ws.Range("N122").Value2 = -15608.2855

$ws.Range("H125").Value2 = 89974.164
$ws.Range("J125").Value2 = 89974.164
$ws.Range("L125").Value2 = 89974.164
$ws.Range("N125").Value2 = -99814.164

$ws.Range("H132").Value2 = 2382.1428
$ws.Range("I132").Value2 = 2382.1428
$ws.Range("K132").Value2 = 7146.428400000001
$ws.Range("M132").Value2 = -4616.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 5847.9
$ws.Range("I132").Value2 = 3826.3142
$ws.Range("K132").Value2 = 11478.9426
$ws.Range("M132").Value2 = -8948.942599999998
